# Update jabatan placeholder to the literal position name coming from Firestore.
$d = $word.ActiveDocument

$d.Content.Find.Execute("{jabatan_orang_1}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Wali Nagari", 2)

# Remove the leading "AN " run in the "AN WALI NAGARI LIMO KOTO" signature line,
# leaving only "WALI NAGARI LIMO KOTO".
$d.Content.Find.Execute("AN WALI NAGARI LIMO KOTO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "WALI NAGARI LIMO KOTO", 2)
